# "Generate Report for Archive"
#
# The localization status report is regenerated: the shared "Status" text
# for the zh-cn / de-de locales moves from "Ready for handoff" to
# "In Translation". This string is surfaced in three places:
#   - Overview!E2 (zh-cn status column) and Overview!F2 (de-de status column)
#   - zh-cn!C2 (Status column)
#   - de-de!C2 (Status column)
#
# Regenerating the report also re-measures the "Status" column to fit the
# (now shorter) text, so that column narrows on every sheet that shows it.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
